$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number and date range) ---
$ws.Range("A8").Value = "Volume 30   Number  14"
$ws.Range("C9").Value = "Report Covering the Week  4/3/2023  Through  4/9/2023"

# --- Cell type/style conversions (copy style+type from a matching donor cell, then correct value) ---
$ws.Range("C23").Copy($ws.Range("D23"))
$ws.Range("D23").Value = 1
$ws.Range("H23").Copy($ws.Range("E23"))
$ws.Range("E23").Value = 0
$ws.Range("D27").Copy($ws.Range("C27"))
$ws.Range("C30").Copy($ws.Range("D30"))
$ws.Range("E22").Copy($ws.Range("E30"))

# --- Plain numeric value updates ---
$ws.Range("F15").Value = 3
$ws.Range("N15").Value = -73.913043478260
$ws.Range("C16").Value = 7
$ws.Range("E16").Value = 133.333333333333
$ws.Range("F16").Value = 23
$ws.Range("G16").Value = 14
$ws.Range("H16").Value = 64.285714285714
$ws.Range("I16").Value = 61
$ws.Range("J16").Value = 38
$ws.Range("K16").Value = 60.526315789473
$ws.Range("L16").Value = 74.285714285714
$ws.Range("M16").Value = -11.594202898550
$ws.Range("N16").Value = -81.230769230769
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 300
$ws.Range("F17").Value = 25
$ws.Range("G17").Value = 17
$ws.Range("H17").Value = 47.058823529411
$ws.Range("I17").Value = 70
$ws.Range("J17").Value = 73
$ws.Range("K17").Value = -4.109589041095
$ws.Range("L17").Value = -2.777777777777
$ws.Range("M17").Value = -11.392405063291
$ws.Range("N17").Value = -65.346534653465
$ws.Range("C18").Value = 5
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 13
$ws.Range("H18").Value = -23.529411764705
$ws.Range("I18").Value = 52
$ws.Range("J18").Value = 60
$ws.Range("K18").Value = -13.333333333333
$ws.Range("L18").Value = 13.043478260869
$ws.Range("M18").Value = 4
$ws.Range("N18").Value = -70.786516853932
$ws.Range("C19").Value = 7
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = 16.666666666666
$ws.Range("F19").Value = 28
$ws.Range("G19").Value = 25
$ws.Range("H19").Value = 12
$ws.Range("I19").Value = 94
$ws.Range("J19").Value = 96
$ws.Range("K19").Value = -2.083333333333
$ws.Range("L19").Value = 1.075268817204
$ws.Range("M19").Value = 25.333333333333
$ws.Range("N19").Value = 1.075268817204
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 8
$ws.Range("G20").Value = 10
$ws.Range("H20").Value = -20
$ws.Range("I20").Value = 28
$ws.Range("J20").Value = 46
$ws.Range("K20").Value = -39.130434782608
$ws.Range("L20").Value = 33.333333333333
$ws.Range("M20").Value = 7.692307692307
$ws.Range("N20").Value = -83.815028901734
$ws.Range("C21").Value = 30
$ws.Range("D21").Value = 19
$ws.Range("E21").Value = 57.894736842105
$ws.Range("F21").Value = 100
$ws.Range("G21").Value = 84
$ws.Range("H21").Value = 19.047619047619
$ws.Range("I21").Value = 312
$ws.Range("J21").Value = 323
$ws.Range("K21").Value = -3.405572755417
$ws.Range("L21").Value = 14.285714285714
$ws.Range("M21").Value = -0.319488817891
$ws.Range("N21").Value = -68.862275449101
$ws.Range("F22").Value = 4
$ws.Range("I22").Value = 6
$ws.Range("K22").Value = 50
$ws.Range("L22").Value = 20
$ws.Range("M22").Value = -14.285714285714
$ws.Range("C23").Value = 1
$ws.Range("F23").Value = 7
$ws.Range("G23").Value = 1
$ws.Range("H23").Value = 600
$ws.Range("I23").Value = 23
$ws.Range("J23").Value = 17
$ws.Range("K23").Value = 35.294117647058
$ws.Range("L23").Value = -17.857142857142
$ws.Range("M23").Value = -4.166666666666
$ws.Range("C24").Value = 30
$ws.Range("D24").Value = 21
$ws.Range("E24").Value = 42.857142857142
$ws.Range("F24").Value = 63
$ws.Range("G24").Value = 61
$ws.Range("H24").Value = 3.278688524590
$ws.Range("I24").Value = 223
$ws.Range("J24").Value = 185
$ws.Range("K24").Value = 20.540540540540
$ws.Range("L24").Value = 70.229007633587
$ws.Range("M24").Value = 29.651162790697
$ws.Range("C25").Value = 10
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = 42.857142857142
$ws.Range("F25").Value = 32
$ws.Range("G25").Value = 28
$ws.Range("H25").Value = 14.285714285714
$ws.Range("I25").Value = 109
$ws.Range("J25").Value = 103
$ws.Range("K25").Value = 5.825242718446
$ws.Range("L25").Value = 51.388888888888
$ws.Range("M25").Value = -39.779005524861
$ws.Range("F26").Value = 6
$ws.Range("H26").Value = 200
$ws.Range("L27").Value = -73.333333333333
$ws.Range("F28").Value = 1
$ws.Range("H28").Value = -50
$ws.Range("L28").Value = 25
$ws.Range("M28").Value = -58.333333333333
$ws.Range("N28").Value = -90.740740740740
$ws.Range("F29").Value = 1
$ws.Range("H29").Value = -50
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -63.636363636363
$ws.Range("N29").Value = -92
